# Updates cryptos list price/volume figures (GitHub Actions style refresh).
#
# Several "Price" cells look numeric (e.g. "1.002", "0.4605") but are stored
# as plain text in the workbook. Assigning such strings straight to
# Range.Value lets Excel's type-inference turn them into real numbers
# (losing the original text representation / introducing float rounding).
# Set-TextValue forces Text format for the duration of the write, then
# clears the formatting again so the cell ends up as a plain, unstyled
# string cell - matching the original file's cell layout.
function Set-TextValue {
    param($Sheet, $CellRef, $Text)
    $Sheet.Range($CellRef).NumberFormat = "@"
    $Sheet.Range($CellRef).Value = $Text
    $Sheet.Range($CellRef).ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "27.581.15"
$ws.Range("E2").Value = "  -1.51%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "1.750.51"
$ws.Range("E3").Value = "  -0.78%  "

# Row 4 - TetherUSD
Set-TextValue $ws "D4" "1.002"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "324.02"
$ws.Range("E5").Value = "  +0.89%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.05%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.4605"
$ws.Range("E7").Value = "  +8.38%  "

# Row 8 - Cardano
Set-TextValue $ws "D8" "0.3588"
$ws.Range("E8").Value = "  -0.46%  "

# Row 9 - Dogecoin
Set-TextValue $ws "D9" "0.07493"
$ws.Range("E9").Value = "  +0.65%  "

# Row 10 - OKB
Set-TextValue $ws "D10" "42.19"
$ws.Range("E10").Value = "  -3.68%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  -0.64%  "

# Row 12 - BinanceUSD
Set-TextValue $ws "D12" "1.002"
$ws.Range("E12").Value = "  -0.11%  "

# Row 13 - Solana
Set-TextValue $ws "D13" "20.70"
$ws.Range("E13").Value = "  -2.21%  "

# Row 14 - Polkadot
Set-TextValue $ws "D14" "5.996"
$ws.Range("E14").Value = "  -1.47%  "

# Row 15 - Chainlink
Set-TextValue $ws "D15" "7.091"
$ws.Range("E15").Value = "  -3.08%  "

# Row 16 - WrappedEther (price only, volume unchanged)
Set-TextValue $ws "D16" "1.751.18"

# Row 17 - Litecoin
Set-TextValue $ws "D17" "92.47"
$ws.Range("E17").Value = "  +1.65%  "

# Row 18 - ShibaInu
Set-TextValue $ws "D18" "0.00001065"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19 - TRON
Set-TextValue $ws "D19" "0.06412"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20 - Dai
Set-TextValue $ws "D20" "1.001"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21 - Avalanche
Set-TextValue $ws "D21" "16.73"
$ws.Range("E21").Value = "  -1.78%  "

# Row 22 - Uniswap
Set-TextValue $ws "D22" "5.812"
$ws.Range("E22").Value = "  -2.61%  "

# Row 23 - WrappedBTC
Set-TextValue $ws "D23" "27.639.73"
$ws.Range("E23").Value = "  -1.32%  "

# Row 24 - Cosmos
Set-TextValue $ws "D24" "11.19"
$ws.Range("E24").Value = "  -0.89%  "

# Row 25 - Toncoin
Set-TextValue $ws "D25" "2.112"
$ws.Range("E25").Value = "  -1.24%  "

# Row 26 - Monero
Set-TextValue $ws "D26" "163.68"
$ws.Range("E26").Value = "  +3.36%  "

# Row 27 - EthereumClassic
Set-TextValue $ws "D27" "20.37"
$ws.Range("E27").Value = "  +1.05%  "

# Row 28 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D28" "1.958.26"
$ws.Range("E28").Value = "  -1.48%  "

# Rows 29/30 swap places: BitcoinCash <-> LidoDAOToken, plus new figures.
# Row 29 becomes LidoDAOToken
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D29" "2.083"
$ws.Range("E29").Value = "  -2.11%  "

# Row 30 becomes BitcoinCash
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D30" "126.78"
$ws.Range("E30").Value = "  +1.30%  "

# Row 31 - ImmutableX
Set-TextValue $ws "D31" "1.073"
$ws.Range("E31").Value = "  -8.01%  "

# Row 32 - Stellar
Set-TextValue $ws "D32" "0.09223"
$ws.Range("E32").Value = "  +3.93%  "

# Row 33 - HuobiToken
Set-TextValue $ws "D33" "3.675"
$ws.Range("E33").Value = "  +3.91%  "

# Row 34 - Filecoin
Set-TextValue $ws "D34" "5.512"
$ws.Range("E34").Value = "  -2.32%  "

# Row 35 - Aptos
Set-TextValue $ws "D35" "11.89"
$ws.Range("E35").Value = "  -5.18%  "

# Row 36 - VeChain
Set-TextValue $ws "D36" "0.02291"
$ws.Range("E36").Value = "  -1.07%  "

# Row 37 - Algorand
Set-TextValue $ws "D37" "0.2100"
$ws.Range("E37").Value = "  -0.26%  "

# Row 38 - Hedera
Set-TextValue $ws "D38" "0.06023"
$ws.Range("E38").Value = "  -0.47%  "

# Row 39 - TheSandbox
Set-TextValue $ws "D39" "0.6340"
$ws.Range("E39").Value = "  -0.49%  "

# Row 40 - InternetComputer(DFINITY)
Set-TextValue $ws "D40" "4.954"
$ws.Range("E40").Value = "  -1.38%  "

# Row 41 - TrustWalletToken
Set-TextValue $ws "D41" "1.197"
$ws.Range("E41").Value = "  +1.11%  "

# Row 42 - WEMIXTOKEN
Set-TextValue $ws "D42" "1.382"
$ws.Range("E42").Value = "  -1.17%  "

# Row 43 - FraxShare
Set-TextValue $ws "D43" "7.767"
$ws.Range("E43").Value = "  -0.89%  "

# Row 44 - EnergySwap
Set-TextValue $ws "D44" "13.18"
$ws.Range("E44").Value = "  -1.91%  "

# Row 45 - Decentraland
Set-TextValue $ws "D45" "0.5893"
$ws.Range("E45").Value = "  -0.46%  "

# Row 46 - PancakeSwap
Set-TextValue $ws "D46" "3.705"
$ws.Range("E46").Value = "  +0.26%  "

# Row 47 - Quant
Set-TextValue $ws "D47" "122.82"
$ws.Range("E47").Value = "  -0.24%  "

# Row 48 - NEARProtocol
Set-TextValue $ws "D48" "1.949"
$ws.Range("E48").Value = "  -2.83%  "

# Row 49 - EOS
Set-TextValue $ws "D49" "1.143"
$ws.Range("E49").Value = "  -4.25%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -0.27%  "

# Row 51 - Aave
Set-TextValue $ws "D51" "72.10"
$ws.Range("E51").Value = "  -2.74%  "
